$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column CO (14-sep) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy header cell CN1 formatting (bold, border, centered) onto CO1
$wsPrix.Range("CN1").Copy()
$wsPrix.Range("CO1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPrix.Range("CO1").Value = "14-sep"
$wsPrix.Range("CO2").Value = 45.73
$wsPrix.Range("CO3").Value = 22.17
$wsPrix.Range("CO4").Value = 17.96
$wsPrix.Range("CO5").Value = 19.06
$wsPrix.Range("CO6").Value = 14.42
$wsPrix.Range("CO7").Value = 16.74
$wsPrix.Range("CO8").Value = 19.38
$wsPrix.Range("CO9").Value = 19.74
$wsPrix.Range("CO10").Value = 18.9
$wsPrix.Range("CO11").Value = 23.58
$wsPrix.Range("CO12").Value = 15.43
$wsPrix.Range("CO13").Value = 17.33
$wsPrix.Range("CO14").Value = 6.5
$wsPrix.Range("CO15").Value = 0
$wsPrix.Range("CO16").Value = -0.01
$wsPrix.Range("CO17").Value = 0
$wsPrix.Range("CO18").Value = 5.59
$wsPrix.Range("CO19").Value = 4.56
$wsPrix.Range("CO20").Value = 15.65
$wsPrix.Range("CO21").Value = 18.4
$wsPrix.Range("CO22").Value = 16.79
$wsPrix.Range("CO23").Value = 13.53
$wsPrix.Range("CO24").Value = 17.36
$wsPrix.Range("CO25").Value = 12.92

# --- Sheet "Gaz": add row 90 (2025-09-12, 32.2) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A90").NumberFormat = "@"
$wsGaz.Range("A90").Value = "2025-09-12"
$wsGaz.Range("A90").ClearFormats()
$wsGaz.Range("B90").Value = 32.2

# --- Sheet "CO2": add row 90 (2025-09-12, 75.47) ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A90").NumberFormat = "@"
$wsCO2.Range("A90").Value = "2025-09-12"
$wsCO2.Range("A90").ClearFormats()
$wsCO2.Range("B90").Value = 75.47

